$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# Row 2
$ws.Range("D2").Value = 171256
$ws.Range("E2").Value = 51095
$ws.Range("F2").Value = 51095
$ws.Range("G2").Value = 50477
$ws.Range("H2").Value = 41952
$ws.Range("I2").Value = 41955
$ws.Range("J2").Value = -3
$ws.Range("K2").Value = 268833
$ws.Range("L2").Value = 88470
$ws.Range("M2").Value = 180363
$ws.Range("N2").Value = 180365
$ws.Range("O2").Value = -2
$ws.Range("P2").Value = 36577
$ws.Range("Q2").Value = 58667
$ws.Range("R2").Value = -60878
$ws.Range("S2").Value = 283
$ws.Range("T2").Value = 48007
$ws.Range("U2").Value = 10660
$ws.Range("V2").Value = 41748
$ws.Range("W2").Value = 29.84
$ws.Range("X2").Value = 24.5
$ws.Range("Y2").Value = 26.98
$ws.Range("Z2").Value = 17.6
$ws.Range("AA2").Value = 49.05
$ws.Range("AB2").Value = 394.26
$ws.Range("AC2").Value = 5842
$ws.Range("AD2").Value = 8.17
$ws.Range("AE2").Value = 24775
$ws.Range("AF2").Value = 1.93
$ws.Range("AG2").Value = 300
$ws.Range("AH2").Value = 0.63
$ws.Range("AI2").Value = 5.21
$ws.Range("AJ2").Value = 728002365

# Row 3
$ws.Range("D3").Value = 187980
$ws.Range("E3").Value = 53361
$ws.Range("F3").Value = 53361
$ws.Range("G3").Value = 52691
$ws.Range("H3").Value = 43236
$ws.Range("I3").Value = 43224
$ws.Range("J3").Value = 12
$ws.Range("K3").Value = 296779
$ws.Range("L3").Value = 82902
$ws.Range("M3").Value = 213877
$ws.Range("N3").Value = 213869
$ws.Range("O3").Value = 8
$ws.Range("P3").Value = 36577
$ws.Range("Q3").Value = 93195
$ws.Range("R3").Value = -71255
$ws.Range("S3").Value = -14623
$ws.Range("T3").Value = 67746
$ws.Range("U3").Value = 25449
$ws.Range("V3").Value = 38186
$ws.Range("W3").Value = 28.39
$ws.Range("X3").Value = 23
$ws.Range("Y3").Value = 21.93
$ws.Range("Z3").Value = 15.29
$ws.Range("AA3").Value = 38.76
$ws.Range("AB3").Value = 505.86
$ws.Range("AC3").Value = 5937
$ws.Range("AD3").Value = 5.18
$ws.Range("AE3").Value = 30293
$ws.Range("AF3").Value = 1.02
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 1.63
$ws.Range("AI3").Value = 8.17
$ws.Range("AJ3").Value = 728002365

# Row 4
$ws.Range("D4").Value = 171980
$ws.Range("E4").Value = 32767
$ws.Range("F4").Value = 32767
$ws.Range("G4").Value = 32165
$ws.Range("H4").Value = 29605
$ws.Range("I4").Value = 29538
$ws.Range("J4").Value = 67
$ws.Range("K4").Value = 322160
$ws.Range("L4").Value = 81925
$ws.Range("M4").Value = 240235
$ws.Range("N4").Value = 240170
$ws.Range("O4").Value = 66
$ws.Range("P4").Value = 36577
$ws.Range("Q4").Value = 55489
$ws.Range("R4").Value = -62305
$ws.Range("S4").Value = 1169
$ws.Range("T4").Value = 59564
$ws.Range("U4").Value = -4074
$ws.Range("V4").Value = 43360
$ws.Range("W4").Value = 19.05
$ws.Range("X4").Value = 17.21
$ws.Range("Y4").Value = 13.01
$ws.Range("Z4").Value = 9.57
$ws.Range("AA4").Value = 34.1
$ws.Range("AB4").Value = 579.89
$ws.Range("AC4").Value = 4057
$ws.Range("AD4").Value = 11.02
$ws.Range("AE4").Value = 34018
$ws.Range("AF4").Value = 1.31
$ws.Range("AG4").Value = 600
$ws.Range("AH4").Value = 1.34
$ws.Range("AI4").Value = 14.34
$ws.Range("AJ4").Value = 728002365

# Row 5
$ws.Range("D5").Value = 301094
$ws.Range("E5").Value = 137213
$ws.Range("F5").Value = 137213
$ws.Range("G5").Value = 134396
$ws.Range("H5").Value = 106422
$ws.Range("I5").Value = 106415
$ws.Range("J5").Value = 7
$ws.Range("K5").Value = 454185
$ws.Range("L5").Value = 115975
$ws.Range("M5").Value = 338209
$ws.Range("N5").Value = 338153
$ws.Range("O5").Value = 56
$ws.Range("P5").Value = 36577
$ws.Range("Q5").Value = 146906
$ws.Range("R5").Value = -119192
$ws.Range("S5").Value = -3519
$ws.Range("T5").Value = 91283
$ws.Range("U5").Value = 55623
$ws.Range("V5").Value = 41713
$ws.Range("W5").Value = 45.57
$ws.Range("X5").Value = 35.34
$ws.Range("Y5").Value = 36.8
$ws.Range("Z5").Value = 27.42
$ws.Range("AA5").Value = 34.29
$ws.Range("AB5").Value = 859.3200000000001
$ws.Range("AC5").Value = 14617
$ws.Range("AD5").Value = 5.23
$ws.Range("AE5").Value = 47897
$ws.Range("AF5").Value = 1.6
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 1.31
$ws.Range("AI5").Value = 6.63
$ws.Range("AJ5").Value = 728002365

# Row 6
$ws.Range("D6").Value = 404451
$ws.Range("E6").Value = 208438
$ws.Range("F6").Value = 208438
$ws.Range("G6").Value = 213410
$ws.Range("H6").Value = 155400
$ws.Range("I6").Value = 155401
$ws.Range("K6").Value = 636583
$ws.Range("L6").Value = 168060
$ws.Range("M6").Value = 468523
$ws.Range("N6").Value = 468457
$ws.Range("P6").Value = 36577
$ws.Range("Q6").Value = 222272
$ws.Range("R6").Value = -214287
$ws.Range("S6").Value = -13953
$ws.Range("T6").Value = 160361
$ws.Range("U6").Value = 61911
$ws.Range("V6").Value = 52819
$ws.Range("W6").Value = 51.54
$ws.Range("X6").Value = 38.42
$ws.Range("Y6").Value = 38.53
$ws.Range("Z6").Value = 28.49
$ws.Range("AA6").Value = 35.87
$ws.Range("AB6").Value = 1262.49
$ws.Range("AC6").Value = 21346
$ws.Range("AD6").Value = 2.83
$ws.Range("AE6").Value = 68488
$ws.Range("AF6").Value = 0.88
$ws.Range("AG6").Value = 1500
$ws.Range("AH6").Value = 2.48
$ws.Range("AI6").Value = 6.6
$ws.Range("AJ6").Value = 728002365

# Row 7
$ws.Range("D7").Value = 268480
$ws.Range("E7").Value = 29317
$ws.Range("G7").Value = 31364
$ws.Range("H7").Value = 24512
$ws.Range("I7").Value = 24493
$ws.Range("K7").Value = 651563
$ws.Range("L7").Value = 168117
$ws.Range("M7").Value = 483446
$ws.Range("N7").Value = 483259
$ws.Range("P7").Value = 36579
$ws.Range("Q7").Value = 100993
$ws.Range("R7").Value = -112970
$ws.Range("S7").Value = 21782
$ws.Range("T7").Value = 132380
$ws.Range("U7").Value = -25639
$ws.Range("W7").Value = 10.92
$ws.Range("X7").Value = 9.130000000000001
$ws.Range("Y7").Value = 5.15
$ws.Range("Z7").Value = 3.81
$ws.Range("AA7").Value = 34.77
$ws.Range("AC7").Value = 3364
$ws.Range("AD7").Value = 29.34
$ws.Range("AE7").Value = 70652
$ws.Range("AF7").Value = 1.4
$ws.Range("AG7").Value = 1158
$ws.Range("AH7").Value = 1.17
$ws.Range("AI7").Value = 34.43

# Row 8
$ws.Range("D8").Value = 315470
$ws.Range("E8").Value = 73743
$ws.Range("G8").Value = 73421
$ws.Range("H8").Value = 56123
$ws.Range("I8").Value = 57101
$ws.Range("K8").Value = 712141
$ws.Range("L8").Value = 179556
$ws.Range("M8").Value = 532585
$ws.Range("N8").Value = 532641
$ws.Range("P8").Value = 36579
$ws.Range("Q8").Value = 148496
$ws.Range("R8").Value = -118528
$ws.Range("S8").Value = -4244
$ws.Range("T8").Value = 105220
$ws.Range("U8").Value = 37917
$ws.Range("W8").Value = 23.38
$ws.Range("X8").Value = 17.79
$ws.Range("Y8").Value = 11.23
$ws.Range("Z8").Value = 8.24
$ws.Range("AA8").Value = 33.71
$ws.Range("AC8").Value = 7844
$ws.Range("AD8").Value = 11.92
$ws.Range("AE8").Value = 77871
$ws.Range("AF8").Value = 1.2
$ws.Range("AG8").Value = 1299
$ws.Range("AH8").Value = 1.39
$ws.Range("AI8").Value = 16.56

# Row 9
$ws.Range("D9").Value = 397081
$ws.Range("E9").Value = 138473
$ws.Range("G9").Value = 141991
$ws.Range("H9").Value = 109148
$ws.Range("I9").Value = 106495
$ws.Range("K9").Value = 821895
$ws.Range("L9").Value = 191152
$ws.Range("M9").Value = 630744
$ws.Range("N9").Value = 632451
$ws.Range("P9").Value = 36579
$ws.Range("Q9").Value = 196781
$ws.Range("R9").Value = -152253
$ws.Range("S9").Value = -13263
$ws.Range("T9").Value = 134454
$ws.Range("U9").Value = 51477
$ws.Range("W9").Value = 34.87
$ws.Range("X9").Value = 27.49
$ws.Range("Y9").Value = 18.28
$ws.Range("Z9").Value = 14.23
$ws.Range("AA9").Value = 30.31
$ws.Range("AC9").Value = 14628
$ws.Range("AD9").Value = 6.39
$ws.Range("AE9").Value = 92463
$ws.Range("AF9").Value = 1.01
$ws.Range("AG9").Value = 1767
$ws.Range("AH9").Value = 1.89
$ws.Range("AI9").Value = 12.08
